$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Header row (row 1): add the missing English column headers (B1:G1 already ---
# --- exist; H1:N1 are brand-new cells that must look like the rest of the header) ---
$ws.Cells.Item(1, 2).Value = "owner"
$ws.Cells.Item(1, 3).Value = "company"
$ws.Cells.Item(1, 4).Value = "address"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "register_date"
$ws.Cells.Item(1, 7).Value = "register_reason"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Match the bold / bordered / centered header look (same xf as B1:G1) for the
# newly added H1:N1 cells.
$headerRange = $ws.Range("H1:N1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- Data row (row 2): fill in the investment / claim / debt / insurance columns ---
$ws.Cells.Item(2, 8).Value = "investment"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "'2011-11-22"
$ws.Cells.Item(2, 11).Value = "王金平"
$ws.Cells.Item(2, 12).Value = 22
$ws.Cells.Item(2, 13).Value = "tmpa85c1"
$ws.Cells.Item(2, 14).Value = 110
